$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 19 ("Buiten Vlaanderen en Brussel", code 99) needs to move
# down to row 20, and a new row 19 ("Niet te lokaliseren", code 93) must
# be inserted in its place, renumbering volgnr accordingly.

$ws.Rows.Item(19).Insert()

# Copy formatting from the row that just moved down (now row 20) onto the
# newly inserted blank row so the new cells keep the same number style.
$ws.Range("A20:B20").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row 19 values.
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 93
$ws.Cells.Item(19, 3).Value = "Niet te lokaliseren"
$ws.Cells.Item(19, 4).Value = "Niet te lokaliseren"

# Renumber the volgnr of the row that shifted down.
$ws.Cells.Item(20, 1).Value = 19
